# Update "想去人数" (F column) values on several sheets to reflect the
# newly generated gh-pages output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 1761
$ws1.Range("F5").Value  = 444
$ws1.Range("F7").Value  = 62
$ws1.Range("F8").Value  = 624
$ws1.Range("F10").Value = 1701
$ws1.Range("F16").Value = 12668
$ws1.Range("F17").Value = 12690
$ws1.Range("F18").Value = 941
$ws1.Range("F21").Value = 500
$ws1.Range("F24").Value = 1987
$ws1.Range("F27").Value = 233
$ws1.Range("F28").Value = 665

# Sheet "演出"
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value  = 13
$ws2.Range("F7").Value  = 6
$ws2.Range("F10").Value = 67
$ws2.Range("F11").Value = 4

# Sheet "本地生活"
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 83
$ws3.Range("F3").Value = 156

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 83
$ws4.Range("F4").Value  = 156
$ws4.Range("F6").Value  = 1761
$ws4.Range("F7").Value  = 444
$ws4.Range("F10").Value = 62
$ws4.Range("F12").Value = 624
$ws4.Range("F15").Value = 1701
$ws4.Range("F22").Value = 12668
$ws4.Range("F23").Value = 12690
$ws4.Range("F24").Value = 941
$ws4.Range("F27").Value = 500
$ws4.Range("F30").Value = 13
$ws4.Range("F31").Value = 6
$ws4.Range("F32").Value = 1987
$ws4.Range("F37").Value = 233
$ws4.Range("F38").Value = 665
$ws4.Range("F39").Value = 67
$ws4.Range("F40").Value = 4

$wb.Save()
